$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: sr no / defect-name / date
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "DashBoard "

$ws.Range("C4").NumberFormat = "d-mmm-yy"
$d = Get-Date -Year 2023 -Month 3 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = $d

# Auto-fit column C to the new date content
$ws.Columns.Item(3).AutoFit()

# Move the active selection to D4, as left by the editor after data entry
$ws.Range("D4").Select()
